# Regenerate merged AHB files
#
# The "L" column ("AENDERUNGSKENNZEICHEN"/change marker) currently shows
# "AENDERUNG" (shared string idx 255, style 7 = bold / goldenrod text) for a
# number of rows. For every such row that sits on a plain (non-highlighted,
# non-NEU, non-ENTFAELLT) segment, the marker is removed entirely: the cell
# becomes blank with the plain centered style (style 4).
#
# In addition, a handful of those rows are the FIRST row of a new Segment
# group (column B changes value) and currently still carry the "white"
# detail-row look (style 5 for most cells / style 7 for the marker cell).
# Those rows are restyled to look like the other group-header rows already
# in the sheet (style 2 everywhere, style 3 for column B, style 4 for the
# now-empty column L) by copying the format of row 2 (an existing, correctly
# styled group-header row) across columns A:V.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are the first row of a new group and need their whole A:V
# formatting turned into the "group header" look (row 2 is a template of
# that look: style 2 on most columns, style 3 on B, style 4 on L).
$headerRows = @(13, 17, 23, 27, 34, 40, 63, 67, 141)

foreach ($rn in $headerRows) {
    $ws.Range("A2:V2").Copy()
    $target = $ws.Range("A" + $rn + ":V" + $rn)
    $target.PasteSpecial(-4122)  # xlPasteFormats
}

# Rows whose "L" cell only loses its "AENDERUNG" marker/style (no other
# column changes). Re-use L2 (style 4, blank) as the format template.
$markerOnlyRows = @(14, 15, 16, 18, 19, 20, 21, 22, 24, 25, 26, 28, 29, 30, 31, 32, 33, 35, 36, 38, 39, 41, 42, 43, 135, 142, 144, 147, 148, 149, 152, 153, 154, 155, 157, 158, 159, 160, 162, 164, 165, 166, 168, 169, 170, 171, 173, 174, 175, 177, 178, 179)

foreach ($rn in $markerOnlyRows) {
    $ws.Range("L2").Copy()
    $target = $ws.Range("L" + $rn)
    $target.PasteSpecial(-4122)  # xlPasteFormats
}

# Finally, blank out the "L" cell's contents for every affected row (the
# format copy above only touches styling, not the cell value).
$allRows = $headerRows + $markerOnlyRows
foreach ($rn in $allRows) {
    $ws.Range("L" + $rn).ClearContents()
}

$excel.CutCopyMode = 0

Write-Output "done"
